$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Backward-extended (earlier) real-time data points to prepend above the
# existing series, which currently starts at 1994-12-31 (row 2).
$newDates = @(
    "1983-12-31",
    "1984-12-31",
    "1985-12-31",
    "1986-12-31",
    "1987-12-31",
    "1988-12-31",
    "1989-12-31",
    "1990-12-31",
    "1991-12-31",
    "1992-12-31",
    "1993-12-31"
)

$newValues = @(
    1.466797881812631,
    2.900424903011278,
    2.603231597845612,
    2.279090113735793,
    1.278816132757377,
    3.44172297297296,
    4.033476219636656,
    5.482086096613448,
    6.118004442050284,
    1.839868480884266,
    -1.202129486518955
)

$shift = $newDates.Length
$oldLastRow = 32
$newLastRow = $oldLastRow + $shift

# Shift the existing data rows (2..32) down by $shift rows, working from the
# bottom up so values are not overwritten before being copied.
for ($r = $oldLastRow; $r -ge 2; $r--) {
    $destRow = $r + $shift
    $ws.Cells.Item($destRow, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($destRow, 2).Value = $ws.Cells.Item($r, 2).Value2
}

# The rows beyond the original extent (33..43) are brand-new cells and did
# not inherit the date column's formatting, so copy it over explicitly.
$ws.Range("A2").Copy()
$ws.Range("A" + ($oldLastRow + 1) + ":A" + $newLastRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the newly freed rows (2..12) with the backward-extended data.
for ($i = 0; $i -lt $shift; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $newDates[$i]
    $ws.Cells.Item($r, 2).Value = $newValues[$i]
}
